# Apply cryptos list update (prices, 1h volume %, and a shift in the
# coin list caused by one coin dropping off the top of the rankings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'60.429.10"
$ws.Range("E2").Value = '  -3.97%  '
$ws.Range("D3").Formula = "'2.904.50"
$ws.Range("E3").Value = '  -3.79%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Formula = "'525.42"
$ws.Range("E5").Value = '  -5.58%  '
$ws.Range("D6").Formula = "'141.35"
$ws.Range("E6").Value = '  -7.52%  '
$ws.Range("D7").Formula = "'0.999"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Formula = "'0.547"
$ws.Range("E8").Value = '  -3.97%  '
$ws.Range("D9").Formula = "'2.908.90"
$ws.Range("E9").Value = '  -3.86%  '
$ws.Range("D10").Formula = "'0.107"
$ws.Range("E10").Value = '  -5.55%  '
$ws.Range("D11").Formula = "'5.88"
$ws.Range("E11").Value = '  -7.25%  '
$ws.Range("D12").Formula = "'0.354"
$ws.Range("E12").Value = '  -3.50%  '
$ws.Range("D13").Formula = "'3.403.71"
$ws.Range("E13").Value = '  -3.99%  '
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("D15").Formula = "'60.505.21"
$ws.Range("E15").Value = '  -3.90%  '
$ws.Range("D16").Formula = "'22.58"
$ws.Range("E16").Value = '  -5.18%  '
$ws.Range("D17").Formula = "'2.899.65"
$ws.Range("E17").Value = '  -3.99%  '
$ws.Range("D18").Formula = "'0.0000139"
$ws.Range("E18").Value = '  -6.84%  '
$ws.Range("D19").Formula = "'4.92"
$ws.Range("E19").Value = '  -3.80%  '
$ws.Range("D20").Formula = "'11.47"
$ws.Range("E20").Value = '  -4.09%  '
$ws.Range("D21").Formula = "'358.91"
$ws.Range("E21").Value = '  -9.17%  '
$ws.Range("D22").Formula = "'6.54"
$ws.Range("E22").Value = '  -2.15%  '
$ws.Range("D23").Formula = "'1.00"
$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Formula = "'63.01"
$ws.Range("E24").Value = '  -3.72%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Formula = "'3.010.49"
$ws.Range("E25").Value = '  -4.36%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Formula = "'0.446"
$ws.Range("E26").Value = '  -4.66%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Formula = "'0.180"
$ws.Range("E27").Value = '  -3.62%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Formula = "'1.00"
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Formula = "'7.78"
$ws.Range("E29").Value = '  -9.47%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Formula = "'0.0₃0848"
$ws.Range("E30").Value = '  -12.92%  '
$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").Formula = "'0.999"
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Formula = "'1.65"
$ws.Range("E32").Value = '  -5.72%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Formula = "'19.35"
$ws.Range("E33").Value = '  -6.11%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Formula = "'150.78"
$ws.Range("E34").Value = '  -6.11%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Formula = "'4.30"
$ws.Range("E35").Value = '  -8.28%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Formula = "'5.52"
$ws.Range("E36").Value = '  -8.84%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Formula = "'0.981"
$ws.Range("E37").Value = '  -10.26%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Formula = "'1.19"
$ws.Range("E38").Value = '  -8.27%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Formula = "'37.96"
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Formula = "'1.47"
$ws.Range("E40").Value = '  -7.13%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Formula = "'2.328.55"
$ws.Range("E41").Value = '  -5.87%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Formula = "'0.644"
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Formula = "'3.63"
$ws.Range("E43").Value = '  -7.66%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Formula = "'20.62"
$ws.Range("E44").Value = '  -8.82%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Formula = "'0.0567"
$ws.Range("E45").Value = '  -5.08%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Formula = "'0.998"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Formula = "'4.85"
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").Formula = "'10.36"
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Formula = "'0.0232"
$ws.Range("E49").Value = '  -6.79%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Formula = "'0.0923"
$ws.Range("E50").Value = '  -3.15%  '
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").Formula = "'248.87"
$ws.Range("E51").Value = '  -5.73%  '
